# Apply the change: remove the "_GoBack" bookmark that wraps the end of the
# SSN paragraph, and add two blank paragraphs right after it (before the
# section properties), giving the findings some breathing room / location.

$d = $word.ActiveDocument

# 1) Remove the leftover "_GoBack" bookmark (bookmarkStart/bookmarkEnd pair)
#    from the last paragraph. "_GoBack" bookmarks are hidden by default, but
#    Bookmarks("_GoBack") / Exists() still reach them directly by name.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) Append two new, empty paragraphs at the very end of the document body
#    (right before the sectPr). Doing this as two separate single "\r"
#    insertions -- each time collapsing to the current end of the story --
#    produces clean empty <w:p/> paragraphs instead of paragraphs that carry
#    a stray empty run.
$end1 = $d.Content
$end1.Collapse(0)
$end1.Text = "`r"

$end2 = $d.Content
$end2.Collapse(0)
$end2.Text = "`r"
